# Commit: "json file change and OS Module"
#
# 1. Sheet1: fill in a 5x multiplication table in columns D:F for rows 2-11
#    (D = 2*n, E = "5*n" text, F = 5*n, for n = 1..10), and move the
#    selection from H14 to J13.
# 2. Add a new Sheet3 at the end of the workbook containing the same
#    multiplication table laid out in columns A:C, and make it the
#    active/selected sheet with the selection on P26.

$wb = $excel.ActiveWorkbook

# --- Sheet1: add D:F multiplication-table columns for rows 2-11 ---
$ws1 = $wb.Worksheets.Item("Sheet1")

for ($n = 1; $n -le 10; $n++) {
    $row = $n + 1
    $ws1.Cells.Item($row, 4).Value = 2 * $n
    $ws1.Cells.Item($row, 5).Value = "5*" + $n
    $ws1.Cells.Item($row, 6).Value = 5 * $n
}

# --- Add Sheet3 at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)

for ($n = 1; $n -le 10; $n++) {
    $ws3.Cells.Item($n, 1).Value = "5*" + $n
    $ws3.Cells.Item($n, 2).Value = 5 * $n
    $ws3.Cells.Item($n, 3).Value = 2 * $n
}

# --- Selections: Sheet1 -> J13, Sheet3 (active) -> P26 ---
[void]$ws1.Range("J13").Select()
[void]$ws3.Activate()
[void]$ws3.Range("P26").Select()
